# Auto-update draw results: append the latest Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 63

# Columns A (date-looking string) and C (digit-only string) would otherwise
# be auto-detected as a date / number by Excel's input parser; a leading
# apostrophe forces them to stay plain text, matching the rest of the sheet.
$ws.Cells.Item($newRow, 1).Value = "'2025-11-18"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "'251118"
$ws.Cells.Item($newRow, 4).Value = "5-3-9-3"
$ws.Cells.Item($newRow, 5).Value = "2025-11-18T21:40:41.879+04:00"
